$d = $word.ActiveDocument

# Mapping of old text -> new text, derived from the unified diff.
$replacements = @(
    @{ Old = "2024-08-24 Saturday"; New = "2024-08-25 Sunday" },
    @{ Old = "125÷4=31, 1";        New = "517÷6=86, 1" },
    @{ Old = "367÷8=45, 7";        New = "766÷2=383, 0" },
    @{ Old = "226÷6=37, 4";        New = "930÷6=155, 0" },
    @{ Old = "397÷3=132, 1";       New = "631÷7=90, 1" },
    @{ Old = "364÷7=52, 0";        New = "623÷3=207, 2" },
    @{ Old = "827÷6=137, 5";       New = "217÷3=72, 1" },
    @{ Old = "695÷4=173, 3";       New = "701÷9=77, 8" },
    @{ Old = "351÷7=50, 1";        New = "825÷7=117, 6" },
    @{ Old = "310÷6=51, 4";        New = "277÷4=69, 1" },
    @{ Old = "660÷2=330, 0";       New = "475÷5=95, 0" },
    @{ Old = "950÷2=475, 0";       New = "462÷3=154, 0" },
    @{ Old = "305÷2=152, 1";       New = "741÷4=185, 1" },
    @{ Old = "792÷7=113, 1";       New = "318÷9=35, 3" },
    @{ Old = "912÷5=182, 2";       New = "325÷5=65, 0" },
    @{ Old = "382÷9=42, 4";        New = "951÷4=237, 3" },
    @{ Old = "524÷5=104, 4";       New = "204÷6=34, 0" },
    @{ Old = "588÷7=84, 0";        New = "443÷3=147, 2" },
    @{ Old = "672÷8=84, 0";        New = "152÷2=76, 0" },
    @{ Old = "993÷4=248, 1";       New = "247÷6=41, 1" },
    @{ Old = "231÷3=77, 0";        New = "176÷2=88, 0" },
    @{ Old = "914÷6=152, 2";       New = "984÷2=492, 0" },
    @{ Old = "759÷3=253, 0";       New = "618÷7=88, 2" },
    @{ Old = "742÷9=82, 4";        New = "697÷8=87, 1" },
    @{ Old = "113÷8=14, 1";        New = "286÷8=35, 6" },
    @{ Old = "539÷6=89, 5";        New = "621÷5=124, 1" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $r.New, 2)
}
